$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E7").Value = 20
$ws.Range("E8").Value = 29
$ws.Range("E14").Value = 30
$ws.Range("E16").Value = 234
$ws.Range("E18").Value = 67
